# "ajout fichier session 2"
# Update the existing planning table and append the new "session 2" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "la page a propos de moi" becomes "la page home"
$ws.Range("C3").Value = "la page home"

# New block for the second work session (rows 8-11)
$ws.Range("B9").Value = "la partie responsive du site"
$ws.Range("C9").Value = "la partie responsive du site"

$ws.Range("B10").Value = "//"
$ws.Range("C10").Value = "//"

$ws.Range("A11").Value = "lundi "
$ws.Range("B11").Value = "mise en forme des bouttons"

$ws.Range("A8").Value = "jeudi "
$ws.Range("B8").Value = "la page login"
$ws.Range("C8").Value = "la page sign up "

$ws.Range("C11").Value = "mise en forme du site avec jquery"

$ws.Range("A9").Value = "vendredi"

# Leave the selection on the last-edited cell, matching the saved view state
[void]$ws.Range("C11").Select()
